$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "index" column data in column A (rows 5-9), mirroring the
# header index row (B1:Y1 contains 0..23) for the data rows of the table.
$ws.Range("A5").Value = 0
$ws.Range("A6").Value = 1
$ws.Range("A7").Value = 2
$ws.Range("A8").Value = 3
$ws.Range("A9").Value = 4

# Widen column A so the new index values are visible (was a near-hidden
# sliver column, now a normal-ish narrow column).
$ws.Columns.Item(1).ColumnWidth = 6.109375

# Update the active selection/cursor position to match the edited area.
[void]$ws.Range("S7:T9").Select()
